$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (was left over near the SSH-key
#    paragraph from the previous edit session).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Find the paragraph that holds "$ ./start.sh" and insert a brand
#    new list paragraph right after it (i.e. right before the
#    "If you want to perform an update ..." list item), carrying the
#    same List Paragraph / numId=8 formatting.
# ------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*`$ ./start.sh*") {
        $targetIndex = $i
        break
    }
}

$nextPara = $d.Paragraphs($targetIndex + 1)
$nextPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($targetIndex + 1)

# Marker trick: type the real sentence plus one throw-away character,
# wrap a bookmark tightly around that throw-away character (so the
# bookmark sits on real, addressable content instead of an ambiguous
# paragraph-boundary point), then delete the throw-away character.
# What remains is the bookmark collapsed exactly at the end of the
# sentence - matching a plain "_GoBack" left by Word after typing.
$newPara.Range.Text = "The default password for the ‘admin’ user is ‘password’. Please change the default password after your first login.X"

$newPara = $d.Paragraphs($targetIndex + 1)
$markRange = $newPara.Range.Duplicate
$markRange.MoveEnd(1, -1) | Out-Null
$markRange.MoveStart(1, ($markRange.End - $markRange.Start) - 1) | Out-Null
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null

$newPara = $d.Paragraphs($targetIndex + 1)
$trailing = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$trailing.Delete() | Out-Null
